$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C (shifts old C.. onward to E..)
$ws.Range("C1:D1").EntireColumn.Insert()

# Fill in the headers for the two newly inserted columns (row 4, like the rest of the header row)
$ws.Range("C4").Value = "FOR ASSIGNING CREATION DATE"
$ws.Range("D4").Value = "ASSIGNMENT CONFIRMATION CREATION DATE"

# The insert pushed the old trailing (empty, formatted-only) O4 cell out to Q4 -- drop it so the
# sheet's used range goes back to ending at column P (no stray empty styled cell).
$ws.Range("Q4").Clear()

# Column width tweaks that came with the new layout (approximate widest achievable value - the
# COM ColumnWidth setter only has ~1/6 character granularity)
$ws.Columns.Item(2).ColumnWidth = 21.307291666666668   # B: STATUS
$ws.Columns.Item(3).ColumnWidth = 51.022135416666664   # C: FOR ASSIGNING CREATION DATE
$ws.Columns.Item(4).ColumnWidth = 46.877604166666664   # D: ASSIGNMENT CONFIRMATION CREATION DATE
$ws.Columns.Item(6).ColumnWidth = 30.877604166666668   # F: SECTION
$ws.Columns.Item(7).ColumnWidth = 26.451822916666668   # G: CATEGORY
$ws.Columns.Item(8).ColumnWidth = 39.736979166666664   # H: NATURE OF PROBLEM
$ws.Columns.Item(10).ColumnWidth = 34.307291666666664  # J: TICKET OWNER

# Selection moves from N4 to A4 and the view scrolls back to show column A
$ws.Range("A4").Select()
